# Apply the "more work towards final product" edit:
#  - carrier ("D") values filled in for the practice rows (2-5)
#  - pair_kind ("J") values filled in for pair rows 6-9 (unique_video / unique_audio)
#  - new rows 14-21 get a "kind" (C) of unique_video / unique_audio and a
#    matching carrier (D) value, mirroring the pattern already used for the
#    generic rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Practice rows: carrier column (D)
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic pair rows: pair_kind column (J)
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New unique_video / unique_audio stimulus rows (kind + carrier)
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
